# The presentation originally carries the "Integral" design (clrScheme
# "Red Violet") on its slide master (ppt/theme/theme1.xml) while the
# notes master keeps the default "Office Theme" (clrScheme "Office",
# ppt/theme/theme2.xml). The edit swaps which theme colours drive the
# slide master / presentation-level theme so it becomes the plain
# "Office Theme" colour scheme (Integral's other settings - fonts,
# fill/line/effect format scheme - are already identical between the
# two themes, only the 12 theme colours differ).
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

# Office theme colour scheme values (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), expressed as the BGR-packed OLE RGB() integers
# PowerPoint's ThemeColorScheme/RGBColor.RGB property uses.
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
